# Add Q4-2022 quarterly data to the 道通科技 (688208) holdings workbook.
#
# Summary of the change:
#   1. A brand-new worksheet "2022-Q4" is inserted right after "总计",
#      pushing every other quarterly sheet one position later.
#   2. The new sheet holds the per-fund holding detail for 2022-Q4.
#   3. The "总计" (totals) sheet gets one new row at the top of its table
#      (right below the header) summarising the new quarter; every
#      previously existing row simply shifts down by one position.

$wb = $excel.ActiveWorkbook

$CENTER = -4108
$TOP = -4160

function Set-HeaderCell($cell, $text) {
    $cell.Value = $text
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = $CENTER
    $cell.VerticalAlignment = $TOP
    $cell.BorderAround(1)
}

function Set-IndexCell($cell, $n) {
    $cell.Value = $n
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = $CENTER
    $cell.VerticalAlignment = $TOP
    $cell.BorderAround(1)
}

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计"
# ---------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2. Fill in the "2022-Q4" per-fund detail sheet
# ---------------------------------------------------------------------

$hB1 = $q4.Range("B1")
Set-HeaderCell $hB1 "基金代码"
$hC1 = $q4.Range("C1")
Set-HeaderCell $hC1 "基金名称"
$hD1 = $q4.Range("D1")
Set-HeaderCell $hD1 "基金规模"
$hE1 = $q4.Range("E1")
Set-HeaderCell $hE1 "股票总仓位"
$hF1 = $q4.Range("F1")
Set-HeaderCell $hF1 "仓位占比"
$hG1 = $q4.Range("G1")
Set-HeaderCell $hG1 "持有市值(亿元)"
$hH1 = $q4.Range("H1")
Set-HeaderCell $hH1 "仓位排名"

# code, name, scale, stockPos, posPct, marketValue, rank
# marketValue = $null means "store literal numeric 0" instead of text "0.00"
$q4rows = @(
    @("470098", "汇添富逆向投资混合A", "16.94", "91.63", "4.86", "0.8233", 6),
    @("470009", "汇添富民营活力混合A", "24.48", "93.48", "2.32", "0.5679", 10),
    @("013680", "华安品质甄选混合A", "12.95", "73.22", "1.80", "0.2331", 6),
    @("001541", "汇添富民营新动力股票", "3.63", "81.10", "2.86", "0.1038", 7),
    @("013681", "华安品质甄选混合C", "5.10", "73.22", "1.80", "0.0918", 6),
    @("002707", "摩根士丹利华鑫科技领先灵活配置混合A", "1.76", "92.23", "4.38", "0.0771", 4),
    @("014509", "汇添富先进制造混合C", "0.93", "86.34", "4.82", "0.0448", 8),
    @("014508", "汇添富先进制造混合A", "0.88", "86.34", "4.82", "0.0424", 8),
    @("015112", "长盛精选行业轮动混合A", "0.17", "53.70", "3.88", "0.0066", 8),
    @("015113", "长盛精选行业轮动混合C", "0.07", "53.70", "3.88", "0.0027", 8),
    @("015182", "汇添富逆向投资混合D", "0.05", "91.63", "4.86", "0.0024", 6),
    @("014871", "摩根士丹利华鑫科技领先灵活配置混合C", "0.05", "92.23", "4.38", "0.0022", 4),
    @("960014", "汇添富民营活力混合 O", "0.00", "93.48", "2.32", $null, 10),
    @("015181", "汇添富逆向投资混合C", "0.00", "91.63", "4.86", $null, 6)
)

$r = 2
foreach ($row in $q4rows) {
    $idx = $r - 2

    $cellA = $q4.Cells.Item($r, 1)
    Set-IndexCell $cellA $idx

    $cellB = $q4.Cells.Item($r, 2)
    Set-TextCell $cellB $row[0]

    $q4.Cells.Item($r, 3).Value = $row[1]

    $cellD = $q4.Cells.Item($r, 4)
    Set-TextCell $cellD $row[2]

    $cellE = $q4.Cells.Item($r, 5)
    Set-TextCell $cellE $row[3]

    $cellF = $q4.Cells.Item($r, 6)
    Set-TextCell $cellF $row[4]

    $mv = $row[5]
    if ($mv -eq $null) {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $cellG = $q4.Cells.Item($r, 7)
        Set-TextCell $cellG $mv
    }

    $q4.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: insert a new 2022-Q4 row and
#    rewrite the whole table (index column + every quarter row) so the
#    0-based index column stays sequential.
# ---------------------------------------------------------------------

$totalRows = @(
    @("2022-Q4", 14, 2),
    @("2022-Q3", 2, 0.09),
    @("2022-Q2", 3, 0.17),
    @("2022-Q1", 6, 1.28),
    @("2021-Q4", 38, 31.08),
    @("2021-Q3", 31, 26.6),
    @("2021-Q2", 38, 25.92),
    @("2021-Q1", 28, 9.449999999999999),
    @("2020-Q4", 9, 5.28)
)

$r = 2
foreach ($row in $totalRows) {
    $idx = $r - 2

    $cellA = $totalSheet.Cells.Item($r, 1)
    Set-IndexCell $cellA $idx

    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]

    $r = $r + 1
}
